$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 8378.84186912827
$ws.Range("C2").Value = 7995.98688358187
$ws.Range("D2").Value = 11764
$ws.Range("E2").Value = 3805.56408513044
$ws.Range("F2").Value = 1.56462369634604

$ws.Range("B3").Value = 8537.0580869407
$ws.Range("C3").Value = 8247.35340606247
$ws.Range("E3").Value = 4139.44194257874
$ws.Range("F3").Value = 140.949806193384

$ws.Range("B4").Value = 8631.53586060117
$ws.Range("C4").Value = 8489.44175752504
$ws.Range("E4").Value = 4648.69027857368
$ws.Range("F4").Value = 172.255501504113

$ws.Range("B5").Value = 8469.55512771084
$ws.Range("C5").Value = 7635.16484036597
$ws.Range("E5").Value = 4583.26303411099
$ws.Range("F5").Value = 133.934494769873

$ws.Range("B6").Value = 2936.96681211716
$ws.Range("C6").Value = 5122.46806117186
$ws.Range("E6").Value = 4120.7351091812
$ws.Range("F6").Value = 9.9667987647108

$ws.Range("B7").Value = 2623.15551157692
$ws.Range("C7").Value = 4911.7107379854
$ws.Range("E7").Value = 3913.99213616859
$ws.Range("F7").Value = -7.42904691025046
